$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data block (rows 2-10, columns A-T) after the edit.
$data = @(
    @("ECs",  "Wnt2b", "Fzd4", "ECs",  3, 1, 2.045113666666667,  6.135341,  0.3607891451960052, 0.3607891451960052, 2, 0.6666666666666666, 13.57958433333333, 40.738753, 0.2289698008477291, 0.2289698008477291, 27.77179350775256, 249.946141569773, 0.08260981872355172, 0.08260981872355172),
    @("ECs",  "Wnt2b", "Fzd4", "FAPs", 3, 1, 2.045113666666667,  6.135341,  0.3607891451960052, 0.3607891451960052, 3, 1,                   19.768727,          59.306181, 0.3333269541315948, 0.3333269541315948, 40.42929376030234, 363.863643842721, 0.1202607468519261,  0.1202607468519261),
    @("ECs",  "Wnt2b", "Fzd4", "sCs",  3, 1, 2.045113666666667,  6.135341,  0.3607891451960052, 0.3607891451960052, 3, 1,                   25.95900466666667, 77.877014, 0.4377032450206762, 0.4377032450206762, 53.08911521686379, 477.802036951774, 0.1579185796205274,  0.1579185796205274),
    @("FAPs", "Wnt2b", "Fzd4", "ECs",  3, 1, 1.870396,            5.611188,  0.3299662923469261, 0.3299662923469261, 2, 0.6666666666666666, 13.57958433333333, 40.738753, 0.2289698008477291, 0.2289698008477291, 25.39920021872934, 228.592801968564, 0.07555231624513922, 0.07555231624513922),
    @("FAPs", "Wnt2b", "Fzd4", "FAPs", 3, 1, 1.870396,            5.611188,  0.3299662923469261, 0.3299662923469261, 3, 1,                   19.768727,          59.306181, 0.3333269541315948, 0.3333269541315948, 36.975347905892,   332.778131153028, 0.1099866591940963,  0.1099866591940963),
    @("FAPs", "Wnt2b", "Fzd4", "sCs",  3, 1, 1.870396,            5.611188,  0.3299662923469261, 0.3299662923469261, 3, 1,                   25.95900466666667, 77.877014, 0.4377032450206762, 0.4377032450206762, 48.55361849251467, 436.982566432632, 0.1444273169076907,  0.1444273169076907),
    @("sCs",  "Wnt2b", "Fzd4", "ECs",  3, 1, 1.752936,            5.258808,  0.3092445624570686, 0.3092445624570686, 2, 0.6666666666666666, 13.57958433333333, 40.738753, 0.2289698008477291, 0.2289698008477291, 23.804142242936,   214.237280186424, 0.07080766587903811, 0.07080766587903811),
    @("sCs",  "Wnt2b", "Fzd4", "FAPs", 3, 1, 1.752936,            5.258808,  0.3092445624570686, 0.3092445624570686, 3, 1,                   19.768727,          59.306181, 0.3333269541315948, 0.3333269541315948, 34.653313232472,   311.879819092248, 0.1030795480855724,  0.1030795480855724),
    @("sCs",  "Wnt2b", "Fzd4", "sCs",  3, 1, 1.752936,            5.258808,  0.3092445624570686, 0.3092445624570686, 3, 1,                   25.95900466666667, 77.877014, 0.4377032450206762, 0.4377032450206762, 45.504473804368,   409.540264239312, 0.1353573484924581,  0.1353573484924581)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = 1 + $j
        $ws.Cells.Item($row, $col).Value = $vals[$j]
    }
}
